# #5: fund, bonds, otherbonds, antique done
# Target sheet: "具有相當價值之財產" (properties of considerable value) = sheet6
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("具有相當價值之財產")

# --- Row 1: replace the (incorrectly duplicated data) header row with real
#     field-name headers, matching the pattern used on the other sheets,
#     and add headers for the newly appended columns F:L ---
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "quantity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "total"
$ws.Cells.Item(1, 6).Value = "property_category"
$ws.Cells.Item(1, 7).Value = "category"
$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"
$ws.Cells.Item(1, 11).Value = "source_file"
$ws.Cells.Item(1, 12).Value = "index"

# give the new header cells (F1:L1) the same style as the existing header
# cells (B1:E1) so they look consistent
$ws.Range("F1:L1").Style = $ws.Range("B1").Style

# --- Data rows 2-13: columns A:E already hold the correct data and are
#     left untouched; append the new F:L columns describing each item as
#     an "otherbonds" property record belonging to the filer ---
$lastRow = 13
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = "otherbonds"
    $ws.Cells.Item($r, 7).Value = "normal"
    $ws.Cells.Item($r, 8).Value = "2013-12-11"
    $ws.Cells.Item($r, 9).Value = "吳育仁"
    $ws.Cells.Item($r, 10).Value = 1734
    $ws.Cells.Item($r, 11).Value = "tmpbcc11"
    $ws.Cells.Item($r, 12).Value = $ws.Cells.Item($r, 1).Value

    # match styling of the existing data columns (B:E use style "2")
    $ws.Range($ws.Cells.Item($r, 6), $ws.Cells.Item($r, 12)).Style = $ws.Range("B2").Style
}

Write-Host "sheet6 (具有相當價值之財產) updated: headers fixed + columns F:L appended"
